# Update the "timestamp" column (Z) on the active sheet of the log workbook.
# Mirrors a re-run of the logging notebook: every previously-logged row gets
# refreshed with the timestamp captured at the moment each chunk of rows was
# written during this run (rows written in the same batch share a timestamp).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ranges = @(
    @{ Rows = "Z2:Z13";  Value = "2025-10-17T07:09:26.729645" },
    @{ Rows = "Z14";     Value = "2025-10-17T07:09:26.743618" },
    @{ Rows = "Z15";     Value = "2025-10-17T07:09:26.744300" },
    @{ Rows = "Z16:Z19"; Value = "2025-10-17T07:09:26.822868" },
    @{ Rows = "Z20:Z23"; Value = "2025-10-17T07:09:26.823867" },
    @{ Rows = "Z24:Z25"; Value = "2025-10-17T07:09:26.824863" },
    @{ Rows = "Z26:Z48"; Value = "2025-10-17T07:09:26.894342" }
)

foreach ($r in $ranges) {
    $ws.Range($r.Rows).Value = $r.Value
}
